# Apply the RBI / MIFOS strategy test-case edits to the loan product workbook.
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoan_Input")
$wsOutput = $wb.Worksheets.Item("ProductLoan_Output")

# --- ProductLoan_Input sheet updates -----------------------------------

# shortname (row 3) switches from the placeholder text "kar5" to the
# numeric product id 391.
$wsInput.Range("B3").Value = 391

# nominalinterestratedefault (row 11) changes from 12 to 1.
$wsInput.Range("B11").Value = 1

# Newly appended ledger-account rows (29-40): key name in column A,
# friendly ledger-account label in column B. Column B (the values) was
# filled in first for all rows, then column A (the keys) - matching the
# shared-string insertion order recorded by the original author.
$newKeys = @(
    "fundsource",
    "loanprotfolio",
    "interestreceivable",
    "penaltiesreceivable",
    "transferinsuspense",
    "feesreceivable",
    "incomefrominterest",
    "incomefrompenalties",
    "incomefromfees",
    "incomefromrecoveryrepayments",
    "loseswrittenoff",
    "overpaymentliability"
)

$newValues = @(
    "Cash",
    "Loan portfolio ",
    "Interest Receivable ",
    "Penalties Receivable ",
    "Transfer in Suspence ",
    "Fees Receivable",
    "Income from interest",
    "Income from penalties",
    "Income from fees",
    "Income from recovery repayments",
    "Losses Writtenoff ",
    "Overpayment Liability"
)

$firstNewRow = 29
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $wsInput.Cells.Item($firstNewRow + $i, 2).Value = $newValues[$i]
}
for ($i = 0; $i -lt $newKeys.Count; $i++) {
    $wsInput.Cells.Item($firstNewRow + $i, 1).Value = $newKeys[$i]
}

# Widen column B to fit the longer ledger-account text now stored there.
$wsInput.Columns.Item(2).ColumnWidth = 61.7109375

# Restore the view so the newly added rows are visible / match the
# recorded selection state.
$wsInput.Application.ActiveWindow.ScrollRow = 13
$wsInput.Range("A68").Select()
